# Add a new test-data row for "GAAP Gl Account Class Profile" to the
# TestData sheet, mirroring the existing rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# New row 7: Iteration number, Description text, SearchText ("Adarsh")
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Enter Text to delete GAAP Gl Account Class Profile"
$ws.Range("C7").Value = "Adarsh"

# Match the updated selection left behind in the sheet view (B7 only,
# instead of the old B7:D10 block).
$ws.Activate() | Out-Null
$ws.Range("B7").Select() | Out-Null
